$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet Folha1 -> Data
$ws.Name = "Data"

# 2) Read the existing vertical (year, value) table before we overwrite it
$years = @()
$vals = @()
for ($r = 2; $r -le 9; $r++) {
    $years += $ws.Cells.Item($r, 1).Value()
    $vals  += $ws.Cells.Item($r, 2).Value()
}
$label = $ws.Cells.Item(1, 2).Value()

# 3) Clear out the old rows 3-9 (no longer needed once transposed)
$ws.Rows("3:9").Delete()

# 4) Write the transposed layout:
#    Row1: A1 = "year" (unchanged label), B1:I1 = the years
#    Row2: A2 = "travels_private" (trimmed label), B2:I2 = the values
for ($i = 0; $i -lt $years.Count; $i++) {
    $col = $i + 2
    $ws.Cells.Item(1, $col).Value = $years[$i]
    $ws.Cells.Item(2, $col).Value = $vals[$i]
}
$ws.Cells.Item(2, 1).Value = "travels_private"

Write-Host "done"
